$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cell B1 from "param" to "name" (ModelWrapper now uses pandas)
$ws.Range("B1").Value = "name"

# Move the active selection to B2
$ws.Range("B2").Select()
